$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.957.06'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '2.596.97'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '523.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.21'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +1.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.71'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.30%  '
$ws.Range('E10').Value = '  +2.19%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('D13').Value = '3.052.63'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('D14').Value = '60.973.86'
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.68'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('E16').Value = '  +1.13%  '
$ws.Range('D17').Value = '2.599.51'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '353.30'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.13'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.427'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.79%  '
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('D26').Value = '2.713.78'
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').Value = '0.0₃0847'
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('E29').Value = '  +0.63%  '
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.31'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +10.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.37'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E33').Value = '  +3.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '148.31'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.20'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.937'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.83%  '
$ws.Range('E37').Value = '  +1.17%  '
$ws.Range('E38').Value = '  +2.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.851'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('E40').Value = '  +1.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.48'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '288.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('E43').Value = '  +1.57%  '
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0561'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('E48').Value = '  +2.38%  '
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.14'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.85%  '
